# Apply metadata updates to the DNA metadata workbook ("added loggers and
# updated metadata"):
#  - Rows 2-97 (column F / "experiment"): tag samples with "PR Rates"
#  - Rows 143-250 (column F / "experiment"): tag samples with "Isotopes",
#    and backfill storage (column I) + sample_size (column K) with the
#    standard values used by the neighbouring rows
#  - Move the active selection to reflect where the author ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-97: only column F (experiment) gets a new "PR Rates" entry.
for ($r = 2; $r -le 97; $r++) {
    $ws.Cells.Item($r, 6).Value = "PR Rates"
}

# Rows 143-250: column F (experiment) gets "Isotopes", column I (storage)
# gets the standard "Molecular -40, shelf 4 right" value, and column K
# (sample_size) gets "0.5 cm".
for ($r = 143; $r -le 250; $r++) {
    $ws.Cells.Item($r, 6).Value = "Isotopes"
    $ws.Cells.Item($r, 9).Value = "Molecular -40, shelf 4 right"
    $ws.Cells.Item($r, 11).Value = "0.5 cm"
}

# Reflect the author's final cursor position/selection.
$ws.Range("F242").Select()
